# Add a new worksheet ("newSheet3") after the last existing sheet,
# then populate it with the "bet type" table (adds an extra "Bet Side"
# column of data, matching the layout already used by the other sheets).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add($null, $lastSheet)
$newWs.Name = "newSheet3"

$newWs.Range("A1").Value = "Board"
$newWs.Range("B1").Value = "Bet Side"
$newWs.Range("C1").Value = "Bet Amt"
$newWs.Range("D1").Value = "Profits"

$newWs.Range("A2").Value = "B"
$newWs.Range("B2").Value = "B"
$newWs.Range("C2").Value = "L2"
$newWs.Range("D2").Value = 2

$newWs.Range("A3").Value = "B"

$newWs.Range("A4").Value = "B"

$newWs.Range("A5").Value = "B"
$newWs.Range("B5").Value = "B"
$newWs.Range("C5").Value = "L3"
$newWs.Range("D5").Value = 5

$newWs.Range("A2").HorizontalAlignment = -4152
$newWs.Range("B2").HorizontalAlignment = -4152
$newWs.Range("A3").HorizontalAlignment = -4152
$newWs.Range("A4").HorizontalAlignment = -4152
$newWs.Range("A5").HorizontalAlignment = -4152
$newWs.Range("B5").HorizontalAlignment = -4152
